$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them (e.g. "522.11" -> 522.11 float)
$textForceCells = @("D5", "D6", "D8", "D9", "D11", "D12", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.681.28'
$ws.Range('E2').Value = '  -2.88%  '
$ws.Range('D3').Value = '2.651.64'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '522.11'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').Value = '143.15'
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').Value = '6.93'
$ws.Range('E9').Value = '  +6.87%  '
$ws.Range('E10').Value = '  -3.97%  '
$ws.Range('D11').Value = '0.334'
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('D12').Value = '0.130'
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('D13').Value = '3.112.84'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').Value = '58.674.60'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').Value = '20.91'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').Value = '2.656.60'
$ws.Range('E16').Value = '  -6.02%  '
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '338.22'
$ws.Range('E18').Value = '  -4.20%  '
$ws.Range('D19').Value = '4.38'
$ws.Range('E19').Value = '  -3.85%  '
$ws.Range('D20').Value = '10.33'
$ws.Range('E20').Value = '  -2.77%  '
$ws.Range('D21').Value = '6.35'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').Value = '64.23'
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('D24').Value = '0.418'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '0.0₃0796'
$ws.Range('E27').Value = '  -3.20%  '
$ws.Range('D28').Value = '7.10'
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('D29').Value = '6.64'
$ws.Range('E29').Value = '  -3.73%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = '1.59'
$ws.Range('E31').Value = '  -0.61%  '
$ws.Range('D32').Value = '18.82'
$ws.Range('E32').Value = '  -1.94%  '
$ws.Range('D33').Value = '150.43'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('D34').Value = '4.12'
$ws.Range('E34').Value = '  -4.61%  '
$ws.Range('D35').Value = '1.18'
$ws.Range('E35').Value = '  -5.37%  '
$ws.Range('D36').Value = '0.903'
$ws.Range('E36').Value = '  -5.47%  '
$ws.Range('D37').Value = '0.860'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').Value = '36.75'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').Value = '1.45'
$ws.Range('E39').Value = '  -5.99%  '
$ws.Range('D40').Value = '3.56'
$ws.Range('E40').Value = '  -3.86%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.33%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '0.611'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').Value = '275.60'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').Value = '19.63'
$ws.Range('E44').Value = '  -2.55%  '
$ws.Range('D45').Value = '0.0966'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').Value = '10.65'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.047.17'
$ws.Range('E47').Value = '  -4.49%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '0.0531'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').Value = '4.69'
$ws.Range('E49').Value = '  -4.90%  '
$ws.Range('D50').Value = '0.0228'
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range('D51').Value = '18.66'
$ws.Range('E51').Value = '  -3.67%  '

# Restore default (General/Normal) style on the forced cells so no stray
# cell-level number formatting is left behind, while keeping the stored
# value as text.
foreach ($c in $textForceCells) {
    $ws.Range($c).Style = "Normal"
}
